# Weekly update: a new Cilantro price record (dated 2022-11-10) is reported
# for the "Macroferia Regional de Talca" market. It becomes the new first
# row of the Cilantro history (row 57), pushing the existing rows 57-79
# down to 58-80 (the former last row, 79, is duplicated down into 80).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 57; Excel shifts rows 57:79 down to 58:80
# and copies the row's formatting (incl. the date-number-format cell D57).
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new weekly record.
$ws.Cells.Item(57, 1).Value2  = 5
$ws.Cells.Item(57, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(57, 3).Value2  = "Maule"
$ws.Cells.Item(57, 4).Value2  = 44875
$ws.Cells.Item(57, 5).Value2  = 7
$ws.Cells.Item(57, 6).Value2  = 100112040
$ws.Cells.Item(57, 7).Value2  = "Cilantro"
$ws.Cells.Item(57, 8).Value2  = "Sin especificar"
$ws.Cells.Item(57, 9).Value2  = "Primera"
$ws.Cells.Item(57, 10).Value2 = 180
$ws.Cells.Item(57, 11).Value2 = 7000
$ws.Cells.Item(57, 12).Value2 = 7000
$ws.Cells.Item(57, 13).Value2 = 7000
$ws.Cells.Item(57, 14).Value2 = "`$/caja 36 atados"
$ws.Cells.Item(57, 15).Value2 = "Región del Maule"
$ws.Cells.Item(57, 16).Value2 = 194
$ws.Cells.Item(57, 17).Value2 = 36
$ws.Cells.Item(57, 18).Value2 = "Hortaliza"
